# 8.5.2.xlsx — add the 2022 column (S) to the "Total unemployment rate"
# table, mirroring the existing 2021 column (R) for values/formatting,
# and move the active selection to T12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new 2022 value for column S (same rows/order as column R = 2021).
# Rows 8 and 36 are section-header rows (no numeric value, just a styled
# blank cell); row 43 repeats the "…" footnote marker used elsewhere in
# that row.
$rowValues = [ordered]@{
    4  = 2022
    5  = 4.9000000000000004
    6  = 6.1
    7  = 4
    9  = 6.1
    10 = 12.4
    11 = 3.2
    12 = 10.8
    13 = 14.6
    14 = 8.5
    15 = 5.5
    16 = 7.1
    17 = 4.4000000000000004
    18 = 5.8
    19 = 11.6
    20 = 3.1
    21 = 1.5
    22 = 2.2999999999999998
    23 = 1
    24 = 2.2999999999999998
    25 = 3.3
    26 = 1.6
    27 = 4.5999999999999996
    28 = 4.4000000000000004
    29 = 4.7
    30 = 4
    31 = 3.2
    32 = 4.7
    33 = 2.6
    34 = 3.3
    35 = 2.2000000000000002
    37 = 13.2
    38 = 7.5
    39 = 4.0999999999999996
    40 = 4.3
    41 = 2.6
    42 = 1
}

# Section-header rows: blank cell, but still needs the R-column formatting
# carried over (bold+italic variant gets created the first time this runs).
$blankRows = @(8, 36)

foreach ($row in $rowValues.Keys) {
    $srcAddr = "R" + $row
    $dstAddr = "S" + $row

    # Copy column R's formatting for this row onto the new S cell, then
    # overwrite with the 2022 figure.
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
    $ws.Range($dstAddr).Value = $rowValues[$row]
}

foreach ($row in $blankRows) {
    $srcAddr = "R" + $row
    $dstAddr = "S" + $row

    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
    $ws.Range($dstAddr).Font.Bold = $true
    $ws.Range($dstAddr).Font.Italic = $true
}

# Row 43: footnote marker "…", same as R43.
$ws.Range("R43").Copy() | Out-Null
$ws.Range("S43").PasteSpecial(-4122) | Out-Null
$ws.Range("S43").Value = $ws.Range("R43").Value()

$excel.CutCopyMode = 0

# Move the live selection, matching the post-edit workbook state.
$ws.Range("T12").Select() | Out-Null
